# Updated symbol list on Tue Feb  7 03:49:24 UTC 2023 with GitHub Actions
# Applies the refreshed Price (D) and Volume(1h) (E) text values for the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "326.74"
Set-TextCell $ws "E2" "0.01%"

Set-TextCell $ws "D3" "44.08"
Set-TextCell $ws "E3" "-2.57%"

Set-TextCell $ws "D4" "5.506"
Set-TextCell $ws "E4" "-1.11%"

Set-TextCell $ws "D5" "0.08025"
Set-TextCell $ws "E5" "-0.82%"

Set-TextCell $ws "D6" "1.996"
Set-TextCell $ws "E6" "4.82%"

Set-TextCell $ws "D8" "0.9487"
Set-TextCell $ws "E8" "-0.09%"

Set-TextCell $ws "D9" "0.1148"
Set-TextCell $ws "E9" "-1.19%"

Set-TextCell $ws "D10" "0.1837"
Set-TextCell $ws "E10" "-3.15%"

Set-TextCell $ws "D11" "12.57"
Set-TextCell $ws "E11" "47.21%"

Set-TextCell $ws "D12" "0.09707"
Set-TextCell $ws "E12" "-3.74%"

Set-TextCell $ws "D13" "0.04599"
Set-TextCell $ws "E13" "9.96%"

Set-TextCell $ws "D14" "0.1067"
Set-TextCell $ws "E14" "0.23%"

Set-TextCell $ws "D15" "0.001253"
Set-TextCell $ws "E15" "-1.75%"

Set-TextCell $ws "D16" "0.04073"
Set-TextCell $ws "E16" "-4.67%"

Set-TextCell $ws "D17" "0.005830"
Set-TextCell $ws "E17" "-3.70%"

Set-TextCell $ws "E18" "-6.91%"

Set-TextCell $ws "D19" "4.285"

Set-TextCell $ws "D20" "0.3478"
Set-TextCell $ws "E20" "-0.24%"

Set-TextCell $ws "D21" "0.1405"
Set-TextCell $ws "E21" "2.27%"

Set-TextCell $ws "D22" "0.2543"
Set-TextCell $ws "E22" "-4.54%"

Set-TextCell $ws "D23" "0.001245"
Set-TextCell $ws "E23" "0.62%"

Set-TextCell $ws "D24" "0.004317"
Set-TextCell $ws "E24" "-6.52%"

Set-TextCell $ws "D25" "0.0001189"
Set-TextCell $ws "E25" "-3.65%"

Set-TextCell $ws "D26" "0.0003740"
Set-TextCell $ws "E26" "-6.53%"

Set-TextCell $ws "D38" "0.02566"
Set-TextCell $ws "E38" "-4.07%"

Set-TextCell $ws "D39" "0.05557"
Set-TextCell $ws "E39" "-0.06%"

Set-TextCell $ws "D40" "0.007547"
Set-TextCell $ws "E40" "-1.98%"

Set-TextCell $ws "D41" "0.1396"

Set-TextCell $ws "D42" "0.007632"
Set-TextCell $ws "E42" "-32.67%"

Set-TextCell $ws "D43" "0.002013"
Set-TextCell $ws "E43" "-3.02%"

Set-TextCell $ws "D44" "0.008521"
Set-TextCell $ws "E44" "-2.13%"

Set-TextCell $ws "D45" "0.00007099"
Set-TextCell $ws "E45" "-0.36%"

Set-TextCell $ws "D46" "0.00000000749"
Set-TextCell $ws "E46" "-0.42%"

Set-TextCell $ws "D47" "0.003527"
Set-TextCell $ws "E47" "54.87%"

Set-TextCell $ws "E48" "-47.10%"

Set-TextCell $ws "D49" "0.00002098"
Set-TextCell $ws "E49" "-0.42%"

Set-TextCell $ws "D50" "0.0001998"
Set-TextCell $ws "E50" "-0.42%"
